$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Price" (column D) values.
# NumberFormat="@" forces text entry (many of these values, e.g. "1.000",
# would otherwise be auto-coerced to a number by Excel); ClearFormats()
# afterwards drops the temporary text format so the cell keeps its original
# (default) style, matching the source data which carries no cell style.
$priceUpdates = @{
    "D2" = "29.909.88"
    "D3" = "1.875.15"
    "D4" = "1.000"
    "D5" = "0.7411"
    "D6" = "242.49"
    "D8" = "0.3149"
    "D9" = "0.07209"
    "D10" = "24.69"
    "D11" = "0.08329"
    "D12" = "0.7503"
    "D13" = "1.897.42"
    "D14" = "5.384"
    "D15" = "92.21"
    "D16" = "6.116"
    "D17" = "29.930.37"
    "D18" = "247.61"
    "D19" = "13.55"
    "D20" = "0.000007846"
    "D21" = "1.000"
    "D22" = "2.124.45"
    "D24" = "1.000"
    "D25" = "0.1551"
    "D26" = "9.297"
    "D27" = "166.09"
    "D30" = "1.492"
    "D31" = "4.578"
    "D33" = "4.225"
    "D34" = "0.05339"
    "D35" = "1.239"
    "D36" = "0.7503"
    "D38" = "2.698"
    "D39" = "0.01962"
    "D40" = "2.755"
    "D41" = "0.4553"
    "D42" = "1.122.73"
    "D43" = "6.132"
    "D44" = "72.46"
    "D45" = "0.8610"
    "D48" = "1.862"
    "D49" = "7.624"
    "D50" = "9.538"
    "D51" = "2.029.15"
}
foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.ClearFormats()
}

# Updated "Volume(1h)" (column E) values - plain text, never numeric-looking
# (they carry surrounding spaces and a trailing "%"), so a direct .Value
# assignment is safe and keeps the inline-string cell type.
$volumeUpdates = @{
    "E2" = "  +0.04%  "
    "E3" = "  -0.74%  "
    "E4" = "  +0.01%  "
    "E5" = "  -3.17%  "
    "E6" = "  -0.12%  "
    "E7" = "  +0.03%  "
    "E8" = "  +0.63%  "
    "E9" = "  +0.56%  "
    "E10" = "  -3.81%  "
    "E11" = "  -2.23%  "
    "E12" = "  -1.58%  "
    "E13" = "  -0.63%  "
    "E14" = "  +0.18%  "
    "E15" = "  -1.72%  "
    "E16" = "  -0.25%  "
    "E17" = "  +0.16%  "
    "E19" = "  -1.65%  "
    "E20" = "  +0.57%  "
    "E21" = "  +0.02%  "
    "E22" = "  +0.37%  "
    "E23" = "  -0.60%  "
    "E24" = "  +0.02%  "
    "E25" = "  -4.21%  "
    "E27" = "  +2.44%  "
    "E28" = "  -0.58%  "
    "E29" = "  -0.60%  "
    "E30" = "  +0.06%  "
    "E31" = "  +1.95%  "
    "E32" = "  +0.29%  "
    "E33" = "  +3.24%  "
    "E34" = "  -1.65%  "
    "E35" = "  -0.44%  "
    "E36" = "  +0.84%  "
    "E37" = "  +0.13%  "
    "E38" = "  +0.19%  "
    "E39" = "  +0.74%  "
    "E40" = "  -0.90%  "
    "E41" = "  +1.93%  "
    "E42" = "  +1.99%  "
    "E43" = "  +1.08%  "
    "E44" = "  -0.76%  "
    "E45" = "  +1.05%  "
    "E46" = "  +1.48%  "
    "E47" = "  +0.12%  "
    "E48" = "  -0.39%  "
    "E49" = "  -0.37%  "
    "E50" = "  -2.44%  "
    "E51" = "  +0.58%  "
}
foreach ($addr in $volumeUpdates.Keys) {
    $ws.Range($addr).Value = $volumeUpdates[$addr]
}
